# BankBranch.xlsx: add Lat/Lon columns (C, D) to the branch table and
# populate every branch row with the same coordinate pair.
# (commit: "excelImport: FieldCastException not resolved")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("C1").Value = "Lat"
$ws.Range("D1").Value = "Lon"

# Populate Lat/Lon for each of the 7 branch rows (rows 2-8)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = 23.8333892
    $ws.Cells.Item($r, 4).Value = 90.414381
}

# Move the selection the way it ended up in the saved workbook
$ws.Range("F8").Select()

# Window chrome moved down the screen between edit sessions
$excel.ActiveWindow.Top = 1305
